$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- MVC code-sample cells (row 7): GameRules -> Gameboard relationship marker ---
$ws.Range("B7").Value = ">"
$ws.Range("C7").Value = ">"
$ws.Range("B7:C7").HorizontalAlignment = -4108  # xlCenter

# --- EventType enum additions: purchase_card / reserve_card get a Card column entry ---
$ws.Range("C11").Value = "purchase_card"
$ws.Range("C13").Value = "reserve_card"

# --- Update the active selection left by the author ---
$ws.Range("H8").Select() | Out-Null
